$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Fill in the previously-empty "R3/R4" column (column 4) for the four
# contributor rows.

$t.Cell(2, 4).Range.Text = "History,`rGeneral bugfixes,`rCmd_history,`rMigration of command_handler`rto process,`rpolling upgrades,`rtesting"

$t.Cell(3, 4).Range.Text = "Help,`rFixed PCB allocation,`rGeneral bugfixes,`rCmd_infinity,`rShort term scheduling with Priority Round Robin,`rtesting"

$t.Cell(4, 4).Range.Text = "Help,`rCmd_alias,`rCmd_alarm,`rGeneral bugfixes,`rUsers manual, Programmers manual,`rtesting"

$t.Cell(5, 4).Range.Text = "Help,`rSys_call,`rInterrupts,`rIrq.s,`rContext switch,`rR3Procs,`rProcess loader,`rGeneral Bugfixes,`rUpdated Search commands,`rFixed PCB allocation,`rtesting"

# The "_GoBack" bookmark that used to sit at the end of row 5's R2 cell
# moves to the end of the text we just typed into row 5's R3/R4 cell.
$d.Bookmarks.Item("_GoBack").Delete()
$t.Cell(5, 4).Range.Bookmarks.Add("_GoBack")

# Resize the table's columns (grid widths only change on this edit).
$t.Columns.Item(1).Width = 67.25
$t.Columns.Item(4).Width = 103.3
$t.Columns.Item(5).Width = 48.6
$t.Columns.Item(6).Width = 48.6
